$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.458.24"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.239.36"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "2.353.58"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "2.579.82"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.52"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "44.261.17"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.14"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.74"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0783"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.87"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.34%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  -6.56%  "
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0294"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.814.39"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +11.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "80.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -3.41%  "
